$d = $word.ActiveDocument
$d.Content.Find.Execute("27.09.2021", $true, $false, $false, $false, $false,
                         $true, 1, $false, "03.10.2022", 2)
